$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '22.402.31'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.572.59'
$ws.Range('E3').Value = '  -3.46%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.000'
$ws.Range('E5').Value = '  -0.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '289.65'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3675'
$ws.Range('E7').Value = '  -2.28%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '49.28'
$ws.Range('E8').Value = '  -1.95%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3394'
$ws.Range('E9').Value = '  -2.30%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.172'
$ws.Range('E10').Value = '  -2.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07641'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.000'
$ws.Range('E12').Value = '  -0.16%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.29'
$ws.Range('E13').Value = '  -2.45%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.076'
$ws.Range('E14').Value = '  -3.41%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.919'
$ws.Range('E15').Value = '  -4.00%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.568.96'
$ws.Range('E16').Value = '  -3.52%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001134'
$ws.Range('E17').Value = '  -4.40%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '89.72'
$ws.Range('E18').Value = '  -5.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06747'
$ws.Range('E19').Value = '  -2.75%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.000'
$ws.Range('E20').Value = '  -0.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.253'
$ws.Range('E21').Value = '  -5.08%  '
$ws.Range('B22').Value = 'Avalanche'
$ws.Range('C22').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '16.55'
$ws.Range('E22').Value = '  -4.02%  '
$ws.Range('B23').Value = 'BitDAO'
$ws.Range('C23').Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.5317'
$ws.Range('E23').Value = '  -7.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.02'
$ws.Range('E24').Value = '  -2.55%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '22.418.48'
$ws.Range('E25').Value = '  -3.91%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.367'
$ws.Range('E26').Value = '  -2.72%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.914'
$ws.Range('E27').Value = '  -2.51%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.02'
$ws.Range('E28').Value = '  -3.72%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '145.89'
$ws.Range('E29').Value = '  -3.18%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.969'
$ws.Range('E30').Value = '  -3.57%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '125.68'
$ws.Range('E31').Value = '  -4.39%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.754.99'
$ws.Range('E32').Value = '  -3.05%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.246'
$ws.Range('E33').Value = '  -6.54%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.021'
$ws.Range('E34').Value = '  +5.10%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.011'
$ws.Range('E35').Value = '  -5.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.17'
$ws.Range('E36').Value = '  -9.48%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.08457'
$ws.Range('E37').Value = '  -3.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02542'
$ws.Range('E38').Value = '  -3.97%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2323'
$ws.Range('E39').Value = '  -3.61%  '
$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.542'
$ws.Range('E40').Value = '  -5.03%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.321'
$ws.Range('E41').Value = '  +2.17%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.06478'
$ws.Range('E42').Value = '  -2.70%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.73'
$ws.Range('E43').Value = '  -7.28%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6366'
$ws.Range('E44').Value = '  -5.94%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.20'
$ws.Range('E45').Value = '  -7.13%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9998'
$ws.Range('E46').Value = '  -0.13%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5995'
$ws.Range('E47').Value = '  -4.74%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.757'
$ws.Range('E48').Value = '  -3.20%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.108'
$ws.Range('E49').Value = '  -5.32%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.259'
$ws.Range('E50').Value = '  +3.63%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '124.93'
$ws.Range('E51').Value = '  -0.87%  '
